$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Odd_Over05_FT / Odd_Under05_FT)
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 15

# Row 4 (Odd_Over05_FT / Odd_Under05_FT / Odd_Over25_FT / Odd_Under25_FT)
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5
